$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$ws.Range("A4").Select()
$win.FreezePanes = $true
$panes = $win.Panes
for ($i=1; $i -le $panes.Count; $i++) {
  $p = $panes.Item($i)
  $members = $p | Get-Member
  Write-Host "PANE $i :"
  Write-Host $members
}
